$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: E1 changes from "Data4" (old position) to "Data4" still,
# but shared-string index shifts because other strings were removed/reordered.
# End result text is the same header row, just E1 stays "Data4".
$ws.Range("E1").Value = "Data4"

# Row 2 data: new "Add Book" API test case values.
$ws.Range("A2").Value = "AddBook"
$ws.Range("B2").Value = "MongoDB"
$ws.Range("C2").Value = 12
$ws.Range("D2").Value = "ty"
$ws.Range("E2").Value = "Ajit Naidu"

# New cell style used by B2:E2 - left horizontal alignment.
$ws.Range("B2:E2").HorizontalAlignment = -4131

# Update the active selection to match the author's final cursor position.
$ws.Range("E8").Select() | Out-Null
